$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column E (Geschwindigkeit), which removes the formulas
# and the header, shifting nothing else left (it's the last column).
$ws.Columns.Item(5).Delete()

# Update the selection to match the target state (E1 column, full column selection)
$ws.Range("E1:E1048576").Select()
